$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was inserted as row 57, pushing the
# existing rows 57-94 down to 58-95.
$ws.Rows(57).Insert()

$ws.Cells.Item(57, 1).Value = 1
$ws.Cells.Item(57, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(57, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(57, 4).Value = 44790
$ws.Cells.Item(57, 5).Value = 15
$ws.Cells.Item(57, 6).Value = 100112021
$ws.Cells.Item(57, 7).Value = "Ají"
$ws.Cells.Item(57, 8).Value = "Inferno"
$ws.Cells.Item(57, 9).Value = "Primera"
$ws.Cells.Item(57, 10).Value = 130
$ws.Cells.Item(57, 11).Value = 9000
$ws.Cells.Item(57, 12).Value = 10000
$ws.Cells.Item(57, 13).Value = 9500
$ws.Cells.Item(57, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(57, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(57, 16).Value = 633
$ws.Cells.Item(57, 17).Value = 15
$ws.Cells.Item(57, 18).Value = "Hortaliza"
